$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D" = -0.0882
    "E" = -0.127
    "G" = 0.4895833333333334
    "H" = 0.4895833333333334
    "I" = 0.7020833333333334
    "J" = 0.7020833333333334
    "K" = 3.44
    "L" = 0.7166666666666667
    "M" = 1.73
    "N" = 0.02255541069100391
    "O" = 0.5029069767441861
    "P" = 1.73
    "Q" = 0.02255541069100391
    "R" = 0.5029069767441861
    "U" = 7.7
    "V" = 0.1003911342894394
    "W" = 0.04231242312423125
    "X" = 0.03846338425444983
    "Y" = 0.00384903886978142
    "Z" = 0.06876790830945559
    "AA" = 0.04828080229226361
    "AB" = 0.03846338425444983
    "AC" = 0.009817418037813788
    "AG" = -7.7
    "AJ" = -0.1115942028985507
    "AK" = -0.0988446726572529
    "AP" = -2.27810650887574
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
